$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: some "Price" values look like plain decimals (single dot) and
# Excel would auto-convert them to numbers on assignment. The source data keeps
# them as text (e.g. "43.085.32" has two dots and can only be text, but "300.07"
# would become the number 300.07). To preserve text typing for those cells we
# temporarily mark the cell as Text before assigning the value, then restore the
# cell style to Normal/General so no stray formatting is left behind.

# --- Row swaps: Chainlink/TRON (rows 12/13), LidoDAOToken/Stellar (rows 41/42) ---

# Swap rows 12 and 13 (Chainlink <-> TRON)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.119"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.05%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.21%  "

# Swap rows 41 and 42 (LidoDAOToken <-> Stellar)
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.58%  "

$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

# --- Price / Volume(1h) updates ---

$ws.Range("D2").Value = "43.063.57"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.308.01"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "2.668.60"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "2.323.43"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").Value = "43.004.05"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.31%  "

$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("E35").Value = "  -3.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.06%  "

$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D43").Value = "2.021.26"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").Value = "2.536.72"
$ws.Range("E51").Value = "  +0.29%  "
